$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("relationships")
$ws.Columns.Item(4).Delete()
